$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 275.375
$ws.Range("I6").Value = 275.375
$ws.Range("K6").Value = 826.125
$ws.Range("M6").Value = -714.125
$ws.Range("H55").Value = 381.75
$ws.Range("I55").Value = 372.625
$ws.Range("J55").Value = 400
$ws.Range("K55").Value = 372.625
$ws.Range("L55").Value = 400
$ws.Range("M55").Value = -158.625
$ws.Range("N55").Value = -828
$ws.Range("H98").Value = 1250.6957
$ws.Range("I98").Value = 1249.6471
$ws.Range("J98").Value = 1253.6666
$ws.Range("K98").Value = 1249.6471
$ws.Range("L98").Value = 1253.6666
$ws.Range("M98").Value = 248.3529000000001
$ws.Range("N98").Value = -4249.6666
$ws.Range("H122").Value = 1250.6957
$ws.Range("I122").Value = 1249.6471
$ws.Range("J122").Value = 1253.6666
$ws.Range("K122").Value = 3748.9413
$ws.Range("L122").Value = 3760.9998
$ws.Range("M122").Value = -1298.9413
$ws.Range("N122").Value = -8660.9998
$ws.Range("H138").Value = 1785.5952
$ws.Range("I138").Value = 1244.375
$ws.Range("J138").Value = 2507.2222
$ws.Range("K138").Value = 3733.125
$ws.Range("L138").Value = 7521.6666
$ws.Range("M138").Value = 1406.875
$ws.Range("N138").Value = -17801.6666

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 766.6667
$ws.Range("I5").Value = 33.666668
$ws.Range("K5").Value = 33.666668
$ws.Range("M5").Value = 78.333332
$ws.Range("H61").Value = 19233710
$ws.Range("I61").Value = 20002420
$ws.Range("K61").Value = 20002420
$ws.Range("M61").Value = -20002208
$ws.Range("H88").Value = 2332.4666
$ws.Range("I88").Value = 2374
$ws.Range("J88").Value = 2285
$ws.Range("K88").Value = 2374
$ws.Range("L88").Value = 2285
$ws.Range("M88").Value = -1968
$ws.Range("N88").Value = -3097
$ws.Range("H91").Value = 2332.4666
$ws.Range("I91").Value = 2374
$ws.Range("J91").Value = 2285
$ws.Range("K91").Value = 2374
$ws.Range("L91").Value = 2285
$ws.Range("M91").Value = -970
$ws.Range("N91").Value = -5093
$ws.Range("H136").Value = 19233710
$ws.Range("I136").Value = 20002420
$ws.Range("K136").Value = 60007260
$ws.Range("M136").Value = -60004710

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 766.6667
$ws.Range("I4").Value = 33.666668
$ws.Range("K4").Value = 33.666668
$ws.Range("M4").Value = 81.333332
$ws.Range("H86").Value = 13890543
$ws.Range("I86").Value = 1673.1904
$ws.Range("J86").Value = 33334962
$ws.Range("K86").Value = 1673.1904
$ws.Range("L86").Value = 33334962
$ws.Range("M86").Value = -550.1904
$ws.Range("N86").Value = -33337208
$ws.Range("H89").Value = 13890543
$ws.Range("I89").Value = 1673.1904
$ws.Range("J89").Value = 33334962
$ws.Range("K89").Value = 8365.951999999999
$ws.Range("L89").Value = 166674810
$ws.Range("M89").Value = -2749.951999999999
$ws.Range("N89").Value = -166686042
$ws.Range("H134").Value = 5421.625
$ws.Range("I134").Value = 4648.0264
$ws.Range("J134").Value = 8361.299999999999
$ws.Range("K134").Value = 13944.0792
$ws.Range("L134").Value = 25083.9
$ws.Range("M134").Value = -11409.0792
$ws.Range("N134").Value = -30153.9

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2292.8572
$ws.Range("I62").Value = 2292.8572
$ws.Range("K62").Value = 2292.8572
$ws.Range("M62").Value = -1668.8572
$ws.Range("H65").Value = 2292.8572
$ws.Range("I65").Value = 2292.8572
$ws.Range("K65").Value = 11464.286
$ws.Range("M65").Value = -8344.286
$ws.Range("H107").Value = 561.36365
$ws.Range("I107").Value = 499.8
$ws.Range("J107").Value = 612.6667
$ws.Range("K107").Value = 499.8
$ws.Range("L107").Value = 612.6667
$ws.Range("M107").Value = 1420.2
$ws.Range("N107").Value = -4452.6667

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1210
$ws.Range("J7").Value = 1658
$ws.Range("L7").Value = 4974
$ws.Range("N7").Value = -5198
$ws.Range("H11").Value = 199.13333
$ws.Range("I11").Value = 143.25
$ws.Range("J11").Value = 422.66666
$ws.Range("K11").Value = 429.75
$ws.Range("L11").Value = 1267.99998
$ws.Range("M11").Value = -289.75
$ws.Range("N11").Value = -1547.99998
$ws.Range("H26").Value = 525
$ws.Range("I26").Value = 300
$ws.Range("K26").Value = 900
$ws.Range("M26").Value = -612
$ws.Range("H41").Value = 1798.75
$ws.Range("I41").Value = 195
$ws.Range("J41").Value = 2333.3333
$ws.Range("K41").Value = 585
$ws.Range("L41").Value = 6999.999899999999
$ws.Range("M41").Value = -247
$ws.Range("N41").Value = -7675.999899999999
$ws.Range("H44").Value = 1452.4783
$ws.Range("I44").Value = 247.16667
$ws.Range("J44").Value = 1877.8823
$ws.Range("K44").Value = 741.50001
$ws.Range("L44").Value = 5633.6469
$ws.Range("M44").Value = -343.50001
$ws.Range("N44").Value = -6429.6469
$ws.Range("H48").Value = 1207.6923
$ws.Range("J48").Value = 1207.6923
$ws.Range("L48").Value = 3623.0769
$ws.Range("N48").Value = -4123.0769
$ws.Range("H88").Value = 1099.3334
$ws.Range("J88").Value = 1170.5714
$ws.Range("L88").Value = 3511.7142
$ws.Range("N88").Value = -4367.7142
$ws.Range("H91").Value = 1099.3334
$ws.Range("J91").Value = 1170.5714
$ws.Range("L91").Value = 3511.7142
$ws.Range("N91").Value = -6475.7142
$ws.Range("H92").Value = 1306
$ws.Range("I92").Value = 1300.75
$ws.Range("J92").Value = 1308.625
$ws.Range("K92").Value = 3902.25
$ws.Range("L92").Value = 3925.875
$ws.Range("M92").Value = -2654.25
$ws.Range("N92").Value = -6421.875
$ws.Range("H98").Value = 327.46155
$ws.Range("I98").Value = 409.75
$ws.Range("K98").Value = 1229.25
$ws.Range("M98").Value = 268.75
$ws.Range("H112").Value = 5102.364
$ws.Range("I112").Value = 2900
$ws.Range("K112").Value = 8700
$ws.Range("M112").Value = -7592
$ws.Range("H115").Value = 2508.375
$ws.Range("J115").Value = 4240
$ws.Range("L115").Value = 12720
$ws.Range("N115").Value = -15070
$ws.Range("H121").Value = 1260.7368
$ws.Range("J121").Value = 1477
$ws.Range("L121").Value = 4431
$ws.Range("N121").Value = -7051
$ws.Range("H125").Value = 5420
$ws.Range("J125").Value = 5420
$ws.Range("L125").Value = 16260
$ws.Range("N125").Value = -26100
$ws.Range("H131").Value = 909.5753
$ws.Range("I131").Value = 504.33334
$ws.Range("K131").Value = 1513.00002
$ws.Range("M131").Value = 3526.99998
$ws.Range("H132").Value = 1824
$ws.Range("I132").Value = 1448.5
$ws.Range("J132").Value = 2575
$ws.Range("K132").Value = 13036.5
$ws.Range("L132").Value = 23175
$ws.Range("M132").Value = -10506.5
$ws.Range("N132").Value = -28235
$ws.Range("H138").Value = 2451
$ws.Range("I138").Value = 1645.5555
$ws.Range("J138").Value = 9700
$ws.Range("K138").Value = 4936.666499999999
$ws.Range("L138").Value = 29100
$ws.Range("M138").Value = 203.3335000000006
$ws.Range("N138").Value = -39380

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 23515.936
$ws.Range("I70").Value = 37282.555
$ws.Range("K70").Value = 37282.555
$ws.Range("M70").Value = -37012.555
$ws.Range("H73").Value = 23515.936
$ws.Range("I73").Value = 37282.555
$ws.Range("K73").Value = 37282.555
$ws.Range("M73").Value = -36346.555
$ws.Range("H122").Value = 3031932.8
$ws.Range("I122").Value = 4168011.2
$ws.Range("J122").Value = 2390.3333
$ws.Range("K122").Value = 12504033.6
$ws.Range("L122").Value = 7170.999899999999
$ws.Range("M122").Value = -12501583.6
$ws.Range("N122").Value = -12070.9999
$ws.Range("H132").Value = 6984.7896
$ws.Range("I132").Value = 5598
$ws.Range("K132").Value = 16794
$ws.Range("M132").Value = -14264

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 855.2632
$ws.Range("I22").Value = 784.1667
$ws.Range("J22").Value = 977.1429000000001
$ws.Range("K22").Value = 784.1667
$ws.Range("L22").Value = 977.1429000000001
$ws.Range("M22").Value = -489.1667
$ws.Range("N22").Value = -1567.1429
$ws.Range("H27").Value = 855.2632
$ws.Range("I27").Value = 784.1667
$ws.Range("J27").Value = 977.1429000000001
$ws.Range("K27").Value = 784.1667
$ws.Range("L27").Value = 977.1429000000001
$ws.Range("M27").Value = -677.1667
$ws.Range("N27").Value = -1191.1429
$ws.Range("H32").Value = 24906.5
$ws.Range("I32").Value = 10013
$ws.Range("J32").Value = 39800
$ws.Range("K32").Value = 10013
$ws.Range("L32").Value = 39800
$ws.Range("M32").Value = -9696
$ws.Range("N32").Value = -40434
